# Apply the commit's text edit: append a new run "  creating a new "
# right after the existing " file" run, at the end of the (only)
# paragraph in the document body.
#
# We deliberately build the addition as its own Range/InsertAfter call
# (rather than a Find&Replace that rewrites " file" in place) so that
# Word emits a brand-new <w:r> element - matching the target XML, which
# keeps the original " file" run untouched and appends a sibling run.

$d = $word.ActiveDocument

# Anchor on the very end of the document's main story so the new text
# lands after "2nd file", still inside the same (only) paragraph and
# before the paragraph mark / sectPr.
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)
$rangeStart = $insertionPoint.Start

$insertionPoint.InsertAfter("  creating a new ")

# The newly inserted text currently has no explicit run formatting; the
# original document explicitly carries <w:lang w:val="en-US"/> on every
# run (see the " file" run immediately before it), so mirror that on the
# freshly typed text - this is exactly what Word does when you keep
# typing at the end of a line that is tagged en-US.
$newRun = $d.Range($rangeStart, $insertionPoint.End)
$newRun.LanguageID = "en-US"
